$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.441.84"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").Value = "2.248.57"
$ws.Range("E3").Value = "  -3.77%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'236.10"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("E6").Value = "  -4.46%  "
$ws.Range("D7").Value = "70.04"
$ws.Range("E7").Value = "  -2.80%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -5.79%  "
$ws.Range("D10").Value = "0.0998"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").Value = "59.34"
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").Value = "'36.70"
$ws.Range("E12").Value = "  +13.62%  "
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "6.77"
$ws.Range("E14").Value = "  -5.29%  "
$ws.Range("D15").Value = "2.583.08"
$ws.Range("E15").Value = "  -3.80%  "
$ws.Range("D16").Value = "'15.20"
$ws.Range("E16").Value = "  -5.26%  "
$ws.Range("D17").Value = "0.863"
$ws.Range("E17").Value = "  -3.72%  "
$ws.Range("D18").Value = "2.248.63"
$ws.Range("E18").Value = "  -3.96%  "
$ws.Range("D19").Value = "42.326.34"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("E21").Value = "  -4.49%  "
$ws.Range("D22").Value = "73.53"
$ws.Range("E22").Value = "  -5.65%  "
$ws.Range("D23").Value = "235.16"
$ws.Range("E23").Value = "  -6.35%  "
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "3.68"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").Value = "2.42"
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  -3.06%  "
$ws.Range("D30").Value = "170.83"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("D31").Value = "20.65"
$ws.Range("E31").Value = "  -6.63%  "
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("E33").Value = "  -5.00%  "
$ws.Range("D34").Value = "0.0723"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").Value = "5.36"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").Value = "4.72"
$ws.Range("E36").Value = "  -6.61%  "
$ws.Range("D37").Value = "'3.70"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "22.45"
$ws.Range("E38").Value = "  +19.92%  "
$ws.Range("D39").Value = "'2.30"
$ws.Range("E39").Value = "  -2.71%  "
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("E41").Value = "  -6.54%  "
$ws.Range("D42").Value = "'65.20"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").Value = "9.33"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("E44").Value = "  -12.11%  "
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("D46").Value = "0.193"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").Value = "4.67"
$ws.Range("E47").Value = "  +14.87%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").Value = "10.23"
$ws.Range("E49").Value = "  +9.76%  "
$ws.Range("D50").Value = "1.19"
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("E51").Value = "  -2.44%  "
